$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Translate character names embedded in English dialogue lines (col D) ---
# <リリー> -> <Lily>
$ws.Range("D4").Value2 = "\n<Lily>Ahh~ I already feel my body healing...`nAhh~ Ooh... Ahh...♥"

# <シィナ> -> <Shina>
$ws.Range("D5").Value2 = "\n<Shina>You sound like a middle-aged man-nya."

# <ライム> -> <Lime>
$ws.Range("D6").Value2 = "\n<Lime>Hey hey-`nI wonder how \n[1] is doing since we started`nour escape game?"

# <ライム> -> <Lime> (short variant, trailing newline preserved)
$ws.Range("D7").Value2 = "\n<Lime>Hey hey-`n"

# <リリー> -> <Lily>
$ws.Range("D8").Value2 = "\n<Lily>Mmm... I wonder..."

# --- Add translated-name column D entries for the character-name rows ---
# Row 2: リリー -> Lily (same as the existing C2 translation)
$ws.Range("D2").Value2 = "Lily"

# Row 29: ライム -> Lime (new translated name)
$ws.Range("D29").Value2 = "Lime"

# Row 30: シィナ -> Shina (same as the existing C30 translation)
$ws.Range("D30").Value2 = "Shina"
